$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the price/volume columns keep their original text formatting
# (values like "134.70" or "0.950" would otherwise be auto-converted to numbers).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "58.374.29"
$ws.Range("E2").Value = "  -4.05%  "
$ws.Range("D3").Value = "2.982.98"
$ws.Range("E3").Value = "  -1.37%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "562.75"
$ws.Range("E5").Value = "  -2.83%  "
$ws.Range("D6").Value = "134.70"
$ws.Range("E6").Value = "  +6.00%  "
$ws.Range("E8").Value = "  +3.48%  "
$ws.Range("D9").Value = "2.976.59"
$ws.Range("E9").Value = "  -1.50%  "
$ws.Range("E10").Value = "  -2.43%  "
$ws.Range("D11").Value = "4.89"
$ws.Range("E11").Value = "  -5.15%  "
$ws.Range("E12").Value = "  +2.03%  "
$ws.Range("E13").Value = "  +1.14%  "
$ws.Range("D14").Value = "33.33"
$ws.Range("E14").Value = "  +1.75%  "
$ws.Range("E15").Value = "  +0.86%  "
$ws.Range("D16").Value = "3.478.00"
$ws.Range("E16").Value = "  -1.31%  "
$ws.Range("D17").Value = "6.88"
$ws.Range("E17").Value = "  +6.97%  "
$ws.Range("D18").Value = "2.977.30"
$ws.Range("D19").Value = "58.359.94"
$ws.Range("E19").Value = "  -4.06%  "
$ws.Range("D20").Value = "424.61"
$ws.Range("E20").Value = "  -2.29%  "
$ws.Range("D21").Value = "13.36"
$ws.Range("E21").Value = "  +1.71%  "
$ws.Range("D22").Value = "0.692"
$ws.Range("E22").Value = "  +3.63%  "
$ws.Range("D23").Value = "7.04"
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").Value = "13.23"
$ws.Range("E24").Value = "  +2.91%  "
$ws.Range("D25").Value = "80.17"
$ws.Range("E25").Value = "  +0.71%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").Value = "0.998"
$ws.Range("E27").Value = "  -0.15%  "
$ws.Range("D28").Value = "2.53"
$ws.Range("E28").Value = "  -2.02%  "
$ws.Range("D29").Value = "7.66"
$ws.Range("E29").Value = "  +4.25%  "
$ws.Range("D30").Value = "2.03"
$ws.Range("E30").Value = "  +5.15%  "
$ws.Range("D31").Value = "25.56"
$ws.Range("E31").Value = "  +0.31%  "
$ws.Range("E32").Value = "  -0.42%  "
$ws.Range("D33").Value = "0.0999"
$ws.Range("E33").Value = "  +6.68%  "
$ws.Range("D34").Value = "5.73"
$ws.Range("E34").Value = "  +1.75%  "
$ws.Range("E35").Value = "  -0.87%  "
$ws.Range("D36").Value = "0.950"
$ws.Range("E36").Value = "  -0.82%  "
$ws.Range("D37").Value = "0.0₃0701"
$ws.Range("E37").Value = "  +4.62%  "
$ws.Range("D38").Value = "48.75"
$ws.Range("E38").Value = "  -4.43%  "
$ws.Range("D39").Value = "8.76"
$ws.Range("E39").Value = "  +3.31%  "
$ws.Range("D40").Value = "2.61"
$ws.Range("E40").Value = "  +4.18%  "
$ws.Range("D41").Value = "0.0353"
$ws.Range("E41").Value = "  -2.25%  "
$ws.Range("E42").Value = "  -1.14%  "
$ws.Range("D43").Value = "382.54"
$ws.Range("E43").Value = "  -1.60%  "
$ws.Range("D44").Value = "2.728.37"
$ws.Range("E44").Value = "  +2.16%  "
$ws.Range("E46").Value = "  +2.96%  "
$ws.Range("D47").Value = "123.27"
$ws.Range("E47").Value = "  +3.29%  "
$ws.Range("E48").Value = "  +3.11%  "
$ws.Range("E49").Value = "  -0.51%  "
$ws.Range("D50").Value = "23.72"
$ws.Range("E50").Value = "  -0.17%  "
$ws.Range("E51").Value = "  +0.90%  "
